$wb = $excel.ActiveWorkbook

# The "AS_classes" sheet had 6 rows describing catering sub-categories
# (E910A..E960A / HSO_0000049..HSO_0000054) removed. Deleting the rows
# shifts the remaining rows up and Excel automatically drops the now
# unused shared strings + renumbers references.
$ws2 = $wb.Worksheets.Item("AS_classes")
$ws2.Rows("28:33").Delete()

# Reflect the author's final cursor position on the AS_classes sheet.
$ws2.Activate()
$ws2.Range("F7").Select()
